# Updates the Price (D) and Volume(1h) (E) columns for the cryptos
# table on the active sheet, matching the refreshed data pulled by the
# scheduled GitHub Actions job. Values are written as literal text (not
# numbers) so strings like "1.00", "0.999" and "58.212.23" keep their
# exact displayed form, matching how the source data is produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "58.212.23"
    "E2" = "  -4.92%  "
    "D3" = "2.567.68"
    "E3" = "  -4.13%  "
    "E4" = "  -0.04%  "
    "D5" = "507.44"
    "E5" = "  -4.84%  "
    "D6" = "144.79"
    "E6" = "  -7.81%  "
    "E7" = "  +0.07%  "
    "D8" = "0.574"
    "E8" = "  -2.79%  "
    "D9" = "2.576.84"
    "E9" = "  -4.48%  "
    "D10" = "6.28"
    "E10" = "  -4.88%  "
    "E11" = "  -5.70%  "
    "E12" = "  -5.63%  "
    "E13" = "  -0.93%  "
    "D14" = "3.018.26"
    "E14" = "  -4.05%  "
    "D15" = "58.236.83"
    "E15" = "  -4.94%  "
    "D16" = "21.04"
    "E16" = "  -5.38%  "
    "E17" = "  -4.96%  "
    "D18" = "2.574.62"
    "E18" = "  -4.43%  "
    "D19" = "4.53"
    "E19" = "  -5.76%  "
    "D20" = "341.64"
    "E20" = "  -4.36%  "
    "E21" = "  -4.66%  "
    "D22" = "6.04"
    "E22" = "  -5.33%  "
    "D23" = "0.999"
    "E23" = "  -0.14%  "
    "D24" = "60.44"
    "E24" = "  -2.44%  "
    "E25" = "  -3.94%  "
    "D26" = "1.00"
    "E26" = "  +0.23%  "
    "D27" = "2.680.79"
    "E27" = "  -4.21%  "
    "E28" = "  -6.11%  "
    "D29" = "0.0₃0811"
    "E29" = "  -6.80%  "
    "D30" = "6.99"
    "E30" = "  -6.02%  "
    "E31" = "  -0.02%  "
    "E32" = "  -1.68%  "
    "D33" = "18.76"
    "E33" = "  -4.57%  "
    "D34" = "148.83"
    "E34" = "  -0.82%  "
    "D35" = "1.53"
    "E35" = "  -6.11%  "
    "D36" = "0.945"
    "E36" = "  +5.91%  "
    "D37" = "3.96"
    "E37" = "  -5.19%  "
    "E38" = "  -6.90%  "
    "D39" = "0.855"
    "E39" = "  -7.01%  "
    "D40" = "36.01"
    "E40" = "  -2.69%  "
    "D41" = "289.80"
    "E41" = "  -5.95%  "
    "E42" = "  -7.23%  "
    "D43" = "3.57"
    "E43" = "  -6.26%  "
    "D44" = "0.0991"
    "E44" = "  -3.41%  "
    "D45" = "0.996"
    "E45" = "  -0.07%  "
    "E46" = "  -7.00%  "
    "E47" = "  -5.89%  "
    "D48" = "19.20"
    "E48" = "  -7.73%  "
    "D49" = "10.25"
    "E49" = "  -0.90%  "
    "E50" = "  -5.83%  "
    "D51" = "4.59"
    "E51" = "  -7.70%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
